# chapter-6.docx: drop the placeholder "1.1  subtitle" sub-heading and its
# paragraph, rewrite the FirstParagraph placeholder body text, and leave the
# "Human-Centricity..." (chapter-6) and "Bibliography" headings as the two
# remaining sections, per the commit's "hardcoded title names" cleanup.

$d = $word.ActiveDocument

# 1. Remove the "1.1  subtitle" Heading 2 paragraph entirely (heading + text).
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text.TrimEnd() -eq "1.1`tsubtitle") {
        $p.Range.Delete()
    }
}

# 2. Swap the leftover "test" placeholder body copy for the real placeholder
#    text used by the template.
$d.Content.Find.Execute("test", $true, $false, $false, $false, $false, `
    $true, 1, $false, "[Target x words]", 2) | Out-Null
